$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the username cell from the old admin address to the new global admin address
$ws.Range("B2").Value = "globaladmin1@cyclotrondev.com"

# Reflect the new active selection on the sheet
$ws.Range("D5").Select()
